$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" date (column C) from 2023-10-09 (45208) to 2023-10-13 (45212)
# for rows 2 through 10.
for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45208) {
        $cell.Value = 45212
    }
}
